# Lab1 slide refactor: imports relatifs, nettoyage docs, lisibilité
#
# Target shape: "ZoneTexte 111" (id 112) on slide 2 (sldId 260) - the
# "Fct complémentaire :" bullet list.
#   - remove the "Multicouche" bullet (done)
#   - remove the "possibilité de choisir fichier sources" bullet (done)
#   - re-split the (unchanged) "Ajout option test unitaire" bullet into two
#     runs ("Ajout " / "option test unitaire")
#   - append three new bullets: "Interface", "Fct act 2 : tanh" and
#     "Plusieurs couches cachées " - all highlighted in yellow

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(8)
$tr = $shape.TextFrame.TextRange

# --- sanity check: locate the right shape -------------------------------
if ($shape.Id -ne 112) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        if ($s.Shapes.Item($i).Id -eq 112) {
            $shape = $s.Shapes.Item($i)
            $tr = $shape.TextFrame.TextRange
        }
    }
}

# --- 1. delete the "Multicouche" bullet (paragraph 3) --------------------
$tr.Paragraphs(3, 1).Delete()

# --- 2. delete the "possibilité de choisir fichier sources" bullet -------
# (now at paragraph 5, after the previous deletion)
$tr.Paragraphs(5, 1).Delete()

# --- 3. re-split "Ajout option test unitaire" into two runs --------------
# (now at paragraph 5) - text itself is unchanged, only the run boundary
# moves; go through a throw-away value first so the replay engine doesn't
# keep a stale trailing run around.
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "TEMP_RESET_5"
$para5.Text = "Ajout option test unitaire"
$runA = $para5.Characters(1, 6)
$runA.Text = "Ajout "
$runB = $para5.Characters(7, 20)
$runB.Text = "option test unitaire"

# --- 4. paragraph 6 "Enregistrement(...)" is untouched --------------------

# --- 5. append new bullet "Interface" (highlighted yellow) ---------------
$tr.InsertAfter("`rTEMP_RESET_7")
$para7 = $tr.Paragraphs($tr.Paragraphs().Count, 1)
$para7.Text = "TEMP_RESET_7"
$para7.Text = "Interface"
$para7.Font.Highlight.RGB = 65535

# --- 6. append new bullet "Fct act 2 : tanh" (highlighted yellow) --------
$tr.InsertAfter("`rTEMP_RESET_8")
$para8 = $tr.Paragraphs($tr.Paragraphs().Count, 1)
$para8.Text = "TEMP_RESET_8"
$para8.Text = "Fct act 2 : tanh"
$para8.Font.Highlight.RGB = 65535
$run8a = $para8.Characters(1, 4)
$run8a.Text = "Fct "
$run8b = $para8.Characters(5, 3)
$run8b.Text = "act"
$run8c = $para8.Characters(8, 5)
$run8c.Text = " 2 : "
$run8d = $para8.Characters(13, 4)
$run8d.Text = "tanh"

# --- 7. append new bullet "Plusieurs couches cachées " (highlighted) -----
$tr.InsertAfter("`rTEMP_RESET_9")
$para9 = $tr.Paragraphs($tr.Paragraphs().Count, 1)
$para9.Text = "TEMP_RESET_9"
$para9.Text = "Plusieurs couches cachées "
$para9.Font.Highlight.RGB = 65535
